$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 70 (G=12604)
$ws.Range("H70").Value = 2107.2307
$ws.Range("I70").Value = 1699.8
$ws.Range("J70").Value = 2361.875
$ws.Range("K70").Value = 5099.4
$ws.Range("L70").Value = 7085.625
$ws.Range("M70").Value = -4829.4
$ws.Range("N70").Value = -7625.625

# Row 73 (G=12604)
$ws.Range("H73").Value = 2107.2307
$ws.Range("I73").Value = 1699.8
$ws.Range("J73").Value = 2361.875
$ws.Range("K73").Value = 5099.4
$ws.Range("L73").Value = 7085.625
$ws.Range("M73").Value = -4163.4
$ws.Range("N73").Value = -8957.625

# Row 76 (G=12602)
$ws.Range("H76").Value = 3982.8333
$ws.Range("I76").Value = 3779.4
$ws.Range("K76").Value = 3779.4
$ws.Range("M76").Value = -3464.4

# Row 79 (G=12602)
$ws.Range("H79").Value = 3982.8333
$ws.Range("I79").Value = 3779.4
$ws.Range("K79").Value = 3779.4
$ws.Range("M79").Value = -2687.4

# Row 125 (G=36228)
$ws.Range("H125").Value = 3307.125
$ws.Range("I125").Value = 2974
$ws.Range("J125").Value = 3507
$ws.Range("K125").Value = 26766
$ws.Range("L125").Value = 31563
$ws.Range("M125").Value = -24306
$ws.Range("N125").Value = -36483

$ws = $wb.Worksheets.Item("ARM")
# Row 45 (G=27714)
$ws.Range("H45").Value = 4022.625
$ws.Range("J45").Value = 4575
$ws.Range("L45").Value = 4575
$ws.Range("N45").Value = -5329

# Row 97 (G=19941)
$ws.Range("H97").Value = 2386.3333
$ws.Range("I97").Value = 2122.125
$ws.Range("K97").Value = 2122.125
$ws.Range("M97").Value = -1626.125

# Row 132 (G=43997)
$ws.Range("H132").Value = 2529.5264
$ws.Range("I132").Value = 2461.7222
$ws.Range("K132").Value = 7385.1666
$ws.Range("M132").Value = -4855.1666

$ws = $wb.Worksheets.Item("BSM")
# Row 22 (G=5092)
$ws.Range("H22").Value = 459
$ws.Range("I22").Value = 472.25
$ws.Range("K22").Value = 472.25
$ws.Range("M22").Value = -299.25

# Row 80 (G=13747)
$ws.Range("H80").Value = 2245.4
$ws.Range("I80").Value = 453
$ws.Range("J80").Value = 3440.3333
$ws.Range("K80").Value = 453
$ws.Range("L80").Value = 3440.3333
$ws.Range("M80").Value = 545
$ws.Range("N80").Value = -5436.3333

# Row 83 (G=13747)
$ws.Range("H83").Value = 2245.4
$ws.Range("I83").Value = 453
$ws.Range("J83").Value = 3440.3333
$ws.Range("K83").Value = 2265
$ws.Range("L83").Value = 17201.6665
$ws.Range("M83").Value = 2727
$ws.Range("N83").Value = -27185.6665

# Row 94 (G=19939)
$ws.Range("H94").Value = 2274.7144
$ws.Range("I94").Value = 2377.0908
$ws.Range("J94").Value = 1899.3334
$ws.Range("K94").Value = 2377.0908
$ws.Range("L94").Value = 1899.3334
$ws.Range("M94").Value = -1926.0908
$ws.Range("N94").Value = -2801.3334

# Row 107 (G=27706)
$ws.Range("H107").Value = 920
$ws.Range("I107").Value = 724
$ws.Range("J107").Value = 1900
$ws.Range("K107").Value = 724
$ws.Range("L107").Value = 1900
$ws.Range("M107").Value = 1196
$ws.Range("N107").Value = -5740

$ws = $wb.Worksheets.Item("CRP")
# Row 7 (G=5361)
$ws.Range("H7").Value = 2095.3157
$ws.Range("J7").Value = 3842.8572
$ws.Range("L7").Value = 3842.8572
$ws.Range("N7").Value = -4068.8572

# Row 15 (G=2406)
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

# Row 22 (G=5367)
$ws.Range("H22").Value = 277
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

# Row 29 (G=2408)
$ws.Range("H29").Value = 4200
$ws.Range("J29").Value = 4200
$ws.Range("L29").Value = 4200
$ws.Range("N29").Value = -4786

# Row 31 (G=44023)
$ws.Range("H31").Value = 2033.9
$ws.Range("I31").Value = 1491.2858
$ws.Range("J31").Value = 3300
$ws.Range("K31").Value = 1491.2858
$ws.Range("L31").Value = 3300
$ws.Range("M31").Value = -1196.2858
$ws.Range("N31").Value = -3890

# Row 34 (G=44023)
$ws.Range("H34").Value = 2033.9
$ws.Range("I34").Value = 1491.2858
$ws.Range("J34").Value = 3300
$ws.Range("K34").Value = 1491.2858
$ws.Range("L34").Value = 3300
$ws.Range("M34").Value = -1289.2858
$ws.Range("N34").Value = -3704

# Row 62 (G=12580)
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

# Row 65 (G=12580)
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

# Row 108 (G=27087)
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

# Row 122 (G=36196)
$ws.Range("H122").Value = 1154
$ws.Range("I122").Value = 1176
$ws.Range("K122").Value = 3528
$ws.Range("M122").Value = -1078

# Row 132 (G=44019)
$ws.Range("H132").Value = 1736.8064
$ws.Range("I132").Value = 1637.1786
$ws.Range("K132").Value = 4911.5358
$ws.Range("M132").Value = -2381.5358

# Row 134 (G=44020)
$ws.Range("H134").Value = 1870.25
$ws.Range("I134").Value = 1766.2727
$ws.Range("K134").Value = 5298.8181
$ws.Range("M134").Value = -2763.8181

$ws = $wb.Worksheets.Item("CUL")
# Row 101 (G=19820)
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

# Row 131 (G=36060)
$ws.Range("H131").Value = 1240.625
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 1240.625
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 3721.875
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -13801.875

$ws = $wb.Worksheets.Item("GSM")
# Row 70 (G=14146)
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

# Row 73 (G=14146)
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

# Row 80 (G=12521)
$ws.Range("H80").Value = 3522
$ws.Range("I80").Value = 1999.6666
$ws.Range("J80").Value = 4092.875
$ws.Range("K80").Value = 1999.6666
$ws.Range("L80").Value = 4092.875
$ws.Range("M80").Value = -1001.6666
$ws.Range("N80").Value = -6088.875

# Row 83 (G=12521)
$ws.Range("H83").Value = 3522
$ws.Range("I83").Value = 1999.6666
$ws.Range("J83").Value = 4092.875
$ws.Range("K83").Value = 9998.333000000001
$ws.Range("L83").Value = 20464.375
$ws.Range("M83").Value = -5006.333000000001
$ws.Range("N83").Value = -30448.375

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (G=5277)
$ws.Range("H22").Value = 7612.375
$ws.Range("I22").Value = 2699.5
$ws.Range("K22").Value = 2699.5
$ws.Range("M22").Value = -2404.5

# Row 27 (G=5277)
$ws.Range("H27").Value = 7612.375
$ws.Range("I27").Value = 2699.5
$ws.Range("K27").Value = 2699.5
$ws.Range("M27").Value = -2592.5

# Row 61 (G=27740)
$ws.Range("H61").Value = 899.6667
$ws.Range("I61").Value = 899.5
$ws.Range("J61").Value = 900
$ws.Range("K61").Value = 899.5
$ws.Range("L61").Value = 900
$ws.Range("M61").Value = -697.5
$ws.Range("N61").Value = -1304

# Row 113 (G=27740)
$ws.Range("H113").Value = 899.6667
$ws.Range("I113").Value = 899.5
$ws.Range("J113").Value = 900
$ws.Range("K113").Value = 899.5
$ws.Range("L113").Value = 900
$ws.Range("M113").Value = 1270.5
$ws.Range("N113").Value = -5240

$ws = $wb.Worksheets.Item("WVR")
# Row 74 (G=19022)
$ws.Range("H74").Value = 28313.4
$ws.Range("I74").Value = 51784
$ws.Range("J74").Value = 12666.333
$ws.Range("K74").Value = 51784
$ws.Range("L74").Value = 12666.333
$ws.Range("M74").Value = -50848
$ws.Range("N74").Value = -14538.333

# Row 77 (G=19022)
$ws.Range("H77").Value = 28313.4
$ws.Range("I77").Value = 51784
$ws.Range("J77").Value = 12666.333
$ws.Range("K77").Value = 155352
$ws.Range("L77").Value = 37998.999
$ws.Range("M77").Value = -150672
$ws.Range("N77").Value = -47358.999

